# Applies weekly price-update diff for Hortaliza / Perejil sheet (Terminal Hortofruticola Agro Chillan).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = (Get-Date -Year 2022 -Month 7 -Day 12 -Hour 0 -Minute 0 -Second 0)
$ws.Range("L2").Value = 750
$ws.Range("M2").Value = 725
$ws.Range("P2").Value = 725

# Row 3
$ws.Range("D3").Value = (Get-Date -Year 2022 -Month 10 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I3").Value = 'Primera'
$ws.Range("J3").Value = 200
$ws.Range("K3").Value = 700
$ws.Range("L3").Value = 800
$ws.Range("M3").Value = 750
$ws.Range("P3").Value = 750

# Row 4
$ws.Range("D4").Value = (Get-Date -Year 2022 -Month 10 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I4").Value = 'Segunda'
$ws.Range("J4").Value = 150
$ws.Range("K4").Value = 600
$ws.Range("L4").Value = 600
$ws.Range("M4").Value = 600
$ws.Range("P4").Value = 600

# Row 5
$ws.Range("D5").Value = (Get-Date -Year 2022 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J5").Value = 160
$ws.Range("K5").Value = 750
$ws.Range("L5").Value = 850
$ws.Range("M5").Value = 800
$ws.Range("P5").Value = 800

# Row 6
$ws.Range("D6").Value = (Get-Date -Year 2022 -Month 8 -Day 26 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J6").Value = 120
$ws.Range("K6").Value = 650
$ws.Range("L6").Value = 650
$ws.Range("M6").Value = 650
$ws.Range("P6").Value = 650

# Row 7
$ws.Range("D7").Value = (Get-Date -Year 2023 -Month 3 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K7").Value = 1500
$ws.Range("L7").Value = 1500
$ws.Range("M7").Value = 1500
$ws.Range("P7").Value = 1500

# Row 8
$ws.Range("D8").Value = (Get-Date -Year 2022 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J8").Value = 300
$ws.Range("K8").Value = 700
$ws.Range("L8").Value = 800
$ws.Range("M8").Value = 750
$ws.Range("P8").Value = 750

# Row 9
$ws.Range("D9").Value = (Get-Date -Year 2022 -Month 9 -Day 27 -Hour 0 -Minute 0 -Second 0)
$ws.Range("K9").Value = 600
$ws.Range("L9").Value = 600
$ws.Range("M9").Value = 600
$ws.Range("P9").Value = 600

# Row 10
$ws.Range("D10").Value = (Get-Date -Year 2022 -Month 11 -Day 18 -Hour 0 -Minute 0 -Second 0)

# Row 11
$ws.Range("D11").Value = (Get-Date -Year 2022 -Month 11 -Day 18 -Hour 0 -Minute 0 -Second 0)

# Row 14
$ws.Range("D14").Value = (Get-Date -Year 2022 -Month 3 -Day 4 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J14").Value = 120
$ws.Range("K14").Value = 650
$ws.Range("L14").Value = 700
$ws.Range("M14").Value = 675
$ws.Range("P14").Value = 675

# Row 15
$ws.Range("D15").Value = (Get-Date -Year 2022 -Month 11 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I15").Value = 'Primera'
$ws.Range("J15").Value = 400
$ws.Range("K15").Value = 700
$ws.Range("L15").Value = 800
$ws.Range("M15").Value = 750
$ws.Range("P15").Value = 750

# Row 16
$ws.Range("D16").Value = (Get-Date -Year 2022 -Month 11 -Day 17 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I16").Value = 'Segunda'
$ws.Range("K16").Value = 600
$ws.Range("L16").Value = 600
$ws.Range("M16").Value = 600
$ws.Range("P16").Value = 600

# Row 17
$ws.Range("D17").Value = (Get-Date -Year 2022 -Month 7 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J17").Value = 200

# Row 18
$ws.Range("D18").Value = (Get-Date -Year 2022 -Month 7 -Day 19 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J18").Value = 150

# Row 19
$ws.Range("D19").Value = (Get-Date -Year 2022 -Month 9 -Day 14 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J19").Value = 300
$ws.Range("K19").Value = 800
$ws.Range("L19").Value = 900
$ws.Range("M19").Value = 850
$ws.Range("P19").Value = 850

# Row 23
$ws.Range("D23").Value = (Get-Date -Year 2022 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J23").Value = 300
$ws.Range("K23").Value = 700
$ws.Range("L23").Value = 800
$ws.Range("M23").Value = 750
$ws.Range("P23").Value = 750

# Row 24
$ws.Range("D24").Value = (Get-Date -Year 2022 -Month 10 -Day 25 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I24").Value = 'Segunda'
$ws.Range("J24").Value = 200
$ws.Range("K24").Value = 600
$ws.Range("L24").Value = 600
$ws.Range("M24").Value = 600
$ws.Range("P24").Value = 600

# Row 25
$ws.Range("D25").Value = (Get-Date -Year 2022 -Month 8 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I25").Value = 'Primera'
$ws.Range("J25").Value = 240
$ws.Range("K25").Value = 750
$ws.Range("L25").Value = 850
$ws.Range("M25").Value = 800
$ws.Range("P25").Value = 800

# Row 26
$ws.Range("D26").Value = (Get-Date -Year 2022 -Month 8 -Day 24 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I26").Value = 'Segunda'
$ws.Range("J26").Value = 200
$ws.Range("K26").Value = 650
$ws.Range("L26").Value = 650
$ws.Range("M26").Value = 650
$ws.Range("P26").Value = 650

# Row 27
$ws.Range("D27").Value = (Get-Date -Year 2022 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I27").Value = 'Primera'
$ws.Range("J27").Value = 200
$ws.Range("K27").Value = 750
$ws.Range("L27").Value = 850
$ws.Range("M27").Value = 800
$ws.Range("P27").Value = 800

# Row 28
$ws.Range("D28").Value = (Get-Date -Year 2022 -Month 8 -Day 31 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I28").Value = 'Segunda'
$ws.Range("J28").Value = 200
$ws.Range("K28").Value = 650
$ws.Range("M28").Value = 650
$ws.Range("P28").Value = 650

# Row 29
$ws.Range("D29").Value = (Get-Date -Year 2022 -Month 11 -Day 3 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J29").Value = 200

# Row 30
$ws.Range("D30").Value = (Get-Date -Year 2022 -Month 2 -Day 18 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I30").Value = 'Primera'
$ws.Range("J30").Value = 100
$ws.Range("L30").Value = 650
$ws.Range("M30").Value = 625
$ws.Range("P30").Value = 625

# Row 31
$ws.Range("D31").Value = (Get-Date -Year 2022 -Month 2 -Day 16 -Hour 0 -Minute 0 -Second 0)
$ws.Range("J31").Value = 120
$ws.Range("K31").Value = 600
$ws.Range("L31").Value = 650
$ws.Range("M31").Value = 625
$ws.Range("P31").Value = 625

# Row 32
$ws.Range("D32").Value = (Get-Date -Year 2022 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("I32").Value = 'Primera'
$ws.Range("J32").Value = 240
$ws.Range("K32").Value = 700
$ws.Range("L32").Value = 800
$ws.Range("M32").Value = 750
$ws.Range("P32").Value = 750

# Row 33
$ws.Range("A33").Value = 7
$ws.Range("B33").Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Range("C33").Value = 'Ñuble'
$ws.Range("D33").Value = (Get-Date -Year 2022 -Month 10 -Day 5 -Hour 0 -Minute 0 -Second 0)
$ws.Range("D33").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("E33").Value = 16
$ws.Range("F33").Value = 100112044
$ws.Range("G33").Value = 'Perejil'
$ws.Range("H33").Value = 'Sin especificar'
$ws.Range("I33").Value = 'Segunda'
$ws.Range("J33").Value = 200
$ws.Range("K33").Value = 600
$ws.Range("L33").Value = 600
$ws.Range("M33").Value = 600
$ws.Range("N33").Value = '$/atado 0,5 a 1 kilo'
$ws.Range("O33").Value = 'Región del Maule'
$ws.Range("P33").Value = 600
$ws.Range("Q33").Value = 1
$ws.Range("R33").Value = 'Hortaliza'

